# Update ExpenseHistory task status strings in the Tasks sheet (column J = "סטטוס").
#
# This corresponds to the commit:
#   "feat: Preserve ExpenseHistory filters in URL for back navigation"
# which updates the status of several task rows:
#   - Row 7  (Report State Loss / Task 6): now done - filters persisted via URL
#   - Row 22 (Reports scope for manager): now pending product decision
#   - Row 24 (required form fields): now pending product decision
#   - Row 29 (budget year for sub-categories): now pending product decision
#   - Row 33 (currency conversion rules): now pending product decision
#   - Rows 38, 39, 41, 42, 43, 44, 45 (the "Nachash"/external table items):
#       now marked N/A, since that table is an external system not in this code

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J7").Value  = "✅ בוצע - פילטרים נשמרים ב-URL ונשמרים בניווט הלוך/חזור"
$ws.Range("J22").Value = "⏳ מחכה להחלטת מוצר - צריך להגדיר היקף דוחות למנג׳ר"
$ws.Range("J24").Value = "⏳ מחכה להחלטת מוצר - צריך להגדיר שדות חובה בטופס"
$ws.Range("J29").Value = "⏳ מחכה להחלטת מוצר - האם להוסיף שנת תקציב לתתי קטגוריות"
$ws.Range("J33").Value = "⏳ מחכה להחלטת מוצר - צריך להגדיר כללים להמרת מטבעות"
$ws.Range("J38").Value = "N/A - טבלת הנחש היא מערכת חיצונית לא בקוד הזה"
$ws.Range("J39").Value = "N/A - טבלת הנחש היא מערכת חיצונית לא בקוד הזה"
$ws.Range("J41").Value = "N/A - טבלת הנחש היא מערכת חיצונית לא בקוד הזה"
$ws.Range("J42").Value = "N/A - טבלת הנחש היא מערכת חיצונית לא בקוד הזה"
$ws.Range("J43").Value = "N/A - טבלת הנחש היא מערכת חיצונית לא בקוד הזה"
$ws.Range("J44").Value = "N/A - טבלת הנחש היא מערכת חיצונית לא בקוד הזה"
$ws.Range("J45").Value = "N/A - טבלת הנחש היא מערכת חיצונית לא בקוד הזה"
